$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 47; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E ("purpose")
    if ($cell.Text -eq "fullRNASEQ") {
        $cell.Value2 = "fullRNASeq"
    }
}
